# Add a "Save" column (column H) to the s_vals sheet, matching the
# existing header styling (copied from G1) and fill the data rows with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: text "Save" with the same style as the other header cells.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Data rows: new "Save" values, all zero, unstyled like the other data cells.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
